# ---------------------------------------------------------------------------
# Update "collected_stats.xlsx": rename Sheet1 -> Default, add a new
# CacheVariants sheet, refresh a few labels on Default, extend the Hardware
# Utilization table with a second (VBMSE) column set, and populate
# CacheVariants with code-variant / cache-set execution-time comparisons.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename first sheet, add the new sheet right after it ---------------
$default = $wb.Worksheets.Item(1)
$default.Name = "Default"

$variants = $wb.Worksheets.Add($null, $default)
$variants.Name = "CacheVariants"

# ===========================================================================
# 2. "Default" sheet tweaks
# ===========================================================================
$default.Activate()

# Row 15 header relabel: "Cached" -> "Cached (16 sets)", drop the old C15
# sub-label entirely.
$default.Range("B15").Value = "Cached (16 sets)"
$default.Range("C15").ClearContents()

# ===========================================================================
# 3. "CacheVariants" sheet content
# ===========================================================================
$variants.Activate()

# --- First comparison table: cache-set counts (32 / 16 / 8 sets) ----------
$variants.Range("B16").Value = "32 sets"
$variants.Range("C16").Value = "Time"
$variants.Range("D16").Value = "Cycles"
$variants.Range("F16").Value = "16 sets"
$variants.Range("G16").Value = "Time"
$variants.Range("H16").Value = "Cycles"
$variants.Range("J16").Value = "8 sets"
$variants.Range("K16").Value = "Time"
$variants.Range("L16").Value = "Cycles"

$variants.Range("B17").Value = "Fib"
$variants.Range("C17").Value = 11756310
$variants.Range("D17").Formula = '=(C17-$H$13)/$D$12'
$variants.Range("F17").Value = "Fib"
$variants.Range("G17").Value = 11756310
$variants.Range("H17").Formula = '=(G17-$H$13)/$D$12'
$variants.Range("J17").Value = "Fib"
$variants.Range("K17").Value = 11756310
$variants.Range("L17").Formula = '=(K17-$H$13)/$D$12'

$variants.Range("B18").Value = "Pow"
$variants.Range("C18").Value = 21170
$variants.Range("D18").Formula = '=(C18-$H$13)/$D$12'
$variants.Range("F18").Value = "Pow"
$variants.Range("G18").Value = 21170
$variants.Range("H18").Formula = '=(G18-$H$13)/$D$12'
$variants.Range("J18").Value = "Pow"
$variants.Range("K18").Value = 21170
$variants.Range("L18").Formula = '=(K18-$H$13)/$D$12'

$variants.Range("B19").Value = "Sort"
$variants.Range("C19").Value = 322090
$variants.Range("D19").Formula = '=(C19-$H$13)/$D$12'
$variants.Range("F19").Value = "Sort"
$variants.Range("G19").Value = 322090
$variants.Range("H19").Formula = '=(G19-$H$13)/$D$12'
$variants.Range("J19").Value = "Sort"
$variants.Range("K19").Value = 322090
$variants.Range("L19").Formula = '=(K19-$H$13)/$D$12'

$variants.Range("B20").Value = "VBSME"
$variants.Range("C20").Value = 2349010
$variants.Range("D20").Formula = '=(C20-$H$13)/$D$12'
$variants.Range("F20").Value = "VBSME"
$variants.Range("G20").Value = 2616810
$variants.Range("H20").Formula = '=(G20-$H$13)/$D$12'
$variants.Range("J20").Value = "VBSME"
$variants.Range("K20").Value = 10790890
$variants.Range("L20").Formula = '=(K20-$H$13)/$D$12'

# Row 15 title, added after the table header row was typed in.
$variants.Range("B15").Value = "Block size = 4"

# --- Second comparison table: block-size variants (16 sets fixed) --------
$variants.Range("B22").Value = "Number sets = 16"

$variants.Range("B23").Value = "Block size = 8"
$variants.Range("C23").Value = "Time"
$variants.Range("D23").Value = "Cycles"
$variants.Range("F23").Value = "Block size = 4"
$variants.Range("G23").Value = "Time"
$variants.Range("H23").Value = "Cycles"
$variants.Range("J23").Value = "Block size = 2"
$variants.Range("K23").Value = "Time"
$variants.Range("L23").Value = "Cycles"

$variants.Range("B24").Value = "Fib"
$variants.Range("C24").Value = 11733650
$variants.Range("D24").Formula = '=(C24-$H$13)/$D$12'
$variants.Range("F24").Value = "Fib"
$variants.Range("G24").Value = 11756310
$variants.Range("H24").Formula = '=(G24-$H$13)/$D$12'
$variants.Range("J24").Value = "Fib"
$variants.Range("K24").Value = 11801630
$variants.Range("L24").Formula = '=(K24-$H$13)/$D$12'

$variants.Range("B25").Value = "Pow"
$variants.Range("C25").Value = 12930
$variants.Range("D25").Formula = '=(C25-$H$13)/$D$12'
$variants.Range("F25").Value = "Pow"
$variants.Range("G25").Value = 21170
$variants.Range("H25").Formula = '=(G25-$H$13)/$D$12'
$variants.Range("J25").Value = "Pow"
$variants.Range("K25").Value = 37650
$variants.Range("L25").Formula = '=(K25-$H$13)/$D$12'

$variants.Range("B26").Value = "Sort"
$variants.Range("C26").Value = 301490
$variants.Range("D26").Formula = '=(C26-$H$13)/$D$12'
$variants.Range("F26").Value = "Sort"
$variants.Range("G26").Value = 322090
$variants.Range("H26").Formula = '=(G26-$H$13)/$D$12'
$variants.Range("J26").Value = "Sort"
$variants.Range("K26").Value = 369470
$variants.Range("L26").Formula = '=(K26-$H$13)/$D$12'

$variants.Range("B27").Value = "VBSME"
$variants.Range("C27").Value = 2206870
$variants.Range("D27").Formula = '=(C27-$H$13)/$D$12'
$variants.Range("F27").Value = "VBSME"
$variants.Range("G27").Value = 2616810
$variants.Range("H27").Formula = '=(G27-$H$13)/$D$12'
$variants.Range("J27").Value = "VBSME"
$variants.Range("K27").Value = 13952990
$variants.Range("L27").Formula = '=(K27-$H$13)/$D$12'

# --- Header / settings block at the top of the sheet ----------------------
$variants.Range("B4").Value = "Default from last sheet is"
$variants.Range("B5").Value = "ADDR_SIZE=12"
$variants.Range("B6").Value = "CACHE_SETS=16"
$variants.Range("B7").Value = "LINE_WORDS=4"
$variants.Range("B8").Value = "MEM_BLOCKS=256"

$variants.Range("B10").Value = "Execution Times"
$variants.Range("B10").Font.Bold = $true

$variants.Range("C11").Value = "First Rise"
$variants.Range("D11").Value = 10
$variants.Range("E11").Value = "ns"

$variants.Range("C12").Value = "Cycle Period"
$variants.Range("D12").Formula = "=2*10"
$variants.Range("E12").Value = "ns"
$variants.Range("H12").Value = "Adjusted"

$variants.Range("C13").Value = "Initialization Time"
$variants.Range("D13").Value = 20
$variants.Range("E13").Value = "ns"
$variants.Range("F13").Value = 2
$variants.Range("G13").Value = "cycles"
$variants.Range("H13").Value = 30

$variants.Columns.Item(2).ColumnWidth = 18.53125
$variants.Columns.Item(6).ColumnWidth = 18.73046875
$variants.Columns.Item(10).ColumnWidth = 18.33203125

# ===========================================================================
# 4. Back to "Default": extend the Hardware Utilization table with a
#    second (VBMSE) code column, added last in the editing session.
# ===========================================================================
$default.Activate()

$default.Range("H23").Font.Bold = $true

$default.Range("F24").Value = "Code: Power"
$default.Range("I24").Value = "Cache"
$default.Range("J24").Value = "No cache"
$default.Range("L24").Value = "Code: VBMSE"

$default.Range("H25").Value = "Utilization"
$default.Range("I25").Value = "Synthesis"
$default.Range("K25").Value = "Out of"

$default.Range("H26").Value = "LUT"
$default.Range("I26").Value = 14425
$default.Range("J26").Value = 12346
$default.Range("K26").Value = 63400

$default.Range("H27").Value = "FF"
$default.Range("I27").Value = 8823
$default.Range("J27").Value = 496
$default.Range("K27").Value = 126800

$default.Range("H28").Value = "BRAM"
$default.Range("I28").Value = 1
$default.Range("J28").Value = 1
$default.Range("K28").Value = 135

$default.Range("H29").Value = "DSP"
$default.Range("I29").Value = 7
$default.Range("J29").Value = 7
$default.Range("K29").Value = 240

$default.Range("H30").Value = "IO"
$default.Range("I30").Value = 130
$default.Range("J30").Value = 130
$default.Range("K30").Value = 210

$default.Range("H31").Value = "BUFG"
$default.Range("I31").Value = 12
$default.Range("J31").Value = 12
$default.Range("K31").Value = 32

$default.Range("I32").Value = "Implementation"

$default.Range("I33").Value = 14415
$default.Range("J33").Value = 12341
$default.Range("K33").Value = 63400

$default.Range("I34").Value = 8823
$default.Range("J34").Value = 496
$default.Range("K34").Value = 126800

$default.Range("I35").Value = 1
$default.Range("J35").Value = 1
$default.Range("K35").Value = 135

$default.Range("I36").Value = 7
$default.Range("J36").Value = 7
$default.Range("K36").Value = 240

$default.Range("I37").Value = 130
$default.Range("J37").Value = 130
$default.Range("K37").Value = 210

$default.Range("I38").Value = 12
$default.Range("J38").Value = 12
$default.Range("K38").Value = 32

# ===========================================================================
# 5. Window / selection state, matching the saved view.
# ===========================================================================
$variants.Activate()
$win2 = $excel.ActiveWindow
$win2.Zoom = 70
$win2.ScrollRow = 4
$variants.Range("K28").Select()

$default.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 13
$win.ScrollColumn = 1
$default.Range("B22").Select()
